# Update points for phone 79174460 -> 0.00 by appending a new row (17)
# matching the existing sheet's layout: A = phone (text), B = birthday
# (blank text), C = total_points (number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

# Column A: phone number, stored as text (matches source data which keeps
# phone numbers as literal strings). Force text formatting before writing
# so Excel doesn't auto-coerce the numeric-looking string to a number,
# then clear the format so the cell keeps the workbook's default style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "79174460"
$ws.Cells.Item($row, 1).ClearFormats()

# Column B: birthday unknown -> empty (but present) text cell, same as
# the other rows with no birthday on file.
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 2).ClearFormats()

# Column C: total_points reset to 0.
$ws.Cells.Item($row, 3).Value = 0
